$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 2 (weekly update - newest record goes on top,
# all existing records shift down by one row; the sheet grows from 57 to
# 58 data rows).
$ws.Rows.Item(2).Insert()

# Reset formatting on the newly inserted row so it matches the plain
# (unstyled) look used by every other data row, instead of inheriting
# the bold/centered header formatting that Insert() copied down.
$ws.Range("A2:R2").Style = "Normal"

# Populate the new row with this week's price record.
$ws.Cells.Item(2, 1).Value = 2
$ws.Cells.Item(2, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(2, 3).Value = "Coquimbo"
$ws.Cells.Item(2, 4).Value = 44756
$ws.Cells.Item(2, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(2, 5).Value = 4
$ws.Cells.Item(2, 6).Value = 100112026
$ws.Cells.Item(2, 7).Value = "Haba"
$ws.Cells.Item(2, 8).Value = "Sin especificar"
$ws.Cells.Item(2, 9).Value = "Primera"
$ws.Cells.Item(2, 10).Value = 500
$ws.Cells.Item(2, 11).Value = 11000
$ws.Cells.Item(2, 12).Value = 12000
$ws.Cells.Item(2, 13).Value = 11500
$ws.Cells.Item(2, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(2, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(2, 16).Value = 460
$ws.Cells.Item(2, 17).Value = 25
$ws.Cells.Item(2, 18).Value = "Hortaliza"
